$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1926605504587156
$ws.Range("C2").Value = 0.5657492354740061
$ws.Range("J2").Value = 0.01834862385321101
$ws.Range("P2").Value = 0.1529051987767584
$ws.Range("S2").Value = 0.07033639143730887
$ws.Range("B3").Value = 0.01075268817204301
$ws.Range("C3").Value = 0.01075268817204301
$ws.Range("J3").Value = 0.03763440860215054
$ws.Range("P3").Value = 0.7526881720430108
$ws.Range("S3").Value = 0.1881720430107527
$ws.Range("O4").Value = 0.02272727272727273
$ws.Range("P4").Value = 0.6818181818181818
$ws.Range("S4").Value = 0.2954545454545455
$ws.Range("P5").Value = 0.5
$ws.Range("S5").Value = 0.5
$ws.Range("B6").Value = 0.06008583690987124
$ws.Range("D6").Value = 0.008583690987124463
$ws.Range("F6").Value = 0.05150214592274678
$ws.Range("J6").Value = 0.3218884120171674
$ws.Range("O6").Value = 0.01716738197424893
$ws.Range("Q6").Value = 0.1416309012875537
$ws.Range("R6").Value = 0.06866952789699571
$ws.Range("S6").Value = 0.3304721030042919
$ws.Range("B7").Value = 0.08843537414965986
$ws.Range("D7").Value = 0.0272108843537415
$ws.Range("E7").Value = 0.006802721088435374
$ws.Range("F7").Value = 0.06802721088435375
$ws.Range("J7").Value = 0.09523809523809523
$ws.Range("Q7").Value = 0.2312925170068027
$ws.Range("R7").Value = 0.04081632653061224
$ws.Range("S7").Value = 0.4421768707482993
$ws.Range("B8").Value = 0.1239130434782609
$ws.Range("D8").Value = 0.02391304347826087
$ws.Range("E8").Value = 0.002173913043478261
$ws.Range("F8").Value = 0.06956521739130435
$ws.Range("J8").Value = 0.1021739130434783
$ws.Range("O8").Value = 0.01521739130434783
$ws.Range("Q8").Value = 0.1543478260869565
$ws.Range("R8").Value = 0.08260869565217391
$ws.Range("S8").Value = 0.4260869565217391
$ws.Range("B9").Value = 0.1209302325581395
$ws.Range("D9").Value = 0.03255813953488372
$ws.Range("F9").Value = 0.07906976744186046
$ws.Range("J9").Value = 0.1069767441860465
$ws.Range("O9").Value = 0.02325581395348837
$ws.Range("Q9").Value = 0.1953488372093023
$ws.Range("R9").Value = 0.09767441860465116
$ws.Range("S9").Value = 0.3441860465116279
$ws.Range("B10").Value = 0.1292576419213974
$ws.Range("D10").Value = 0.01834061135371179
$ws.Range("E10").Value = 0.0008733624454148472
$ws.Range("F10").Value = 0.07161572052401746
$ws.Range("J10").Value = 0.1117903930131004
$ws.Range("O10").Value = 0.01572052401746725
$ws.Range("Q10").Value = 0.1947598253275109
$ws.Range("R10").Value = 0.08820960698689956
$ws.Range("S10").Value = 0.3694323144104804
$ws.Range("G11").Value = 0.1266968325791855
$ws.Range("J11").Value = 0.07692307692307693
$ws.Range("K11").Value = 0.1719457013574661
$ws.Range("L11").Value = 0.6244343891402715
$ws.Range("G12").Value = 0.6573426573426573
$ws.Range("J12").Value = 0.2587412587412588
$ws.Range("K12").Value = 0.01398601398601399
$ws.Range("L12").Value = 0.02797202797202797
$ws.Range("S12").Value = 0.04195804195804196
$ws.Range("G13").Value = 0.7391304347826086
$ws.Range("J13").Value = 0.2173913043478261
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("F15").Value = 0.02884615384615385
$ws.Range("H15").Value = 0.1730769230769231
$ws.Range("I15").Value = 0.0673076923076923
$ws.Range("J15").Value = 0.3894230769230769
$ws.Range("K15").Value = 0.03846153846153846
$ws.Range("M15").Value = 0.01923076923076923
$ws.Range("O15").Value = 0.03846153846153846
$ws.Range("S15").Value = 0.2451923076923077
$ws.Range("F16").Value = 0.02816901408450704
$ws.Range("H16").Value = 0.2065727699530517
$ws.Range("I16").Value = 0.1220657276995305
$ws.Range("J16").Value = 0.3802816901408451
$ws.Range("K16").Value = 0.07511737089201878
$ws.Range("M16").Value = 0.02347417840375587
$ws.Range("O16").Value = 0.05164319248826291
$ws.Range("S16").Value = 0.1126760563380282
$ws.Range("F17").Value = 0.03
$ws.Range("H17").Value = 0.1975
$ws.Range("I17").Value = 0.135
$ws.Range("J17").Value = 0.3375
$ws.Range("K17").Value = 0.0925
$ws.Range("M17").Value = 0.0175
$ws.Range("O17").Value = 0.055
$ws.Range("S17").Value = 0.135
$ws.Range("F18").Value = 0.01111111111111111
$ws.Range("H18").Value = 0.1722222222222222
$ws.Range("I18").Value = 0.1
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.08888888888888889
$ws.Range("M18").Value = 0.02777777777777778
$ws.Range("O18").Value = 0.07222222222222222
$ws.Range("S18").Value = 0.1277777777777778
$ws.Range("F19").Value = 0.02544529262086514
$ws.Range("H19").Value = 0.2298558100084818
$ws.Range("I19").Value = 0.08821034775233248
$ws.Range("J19").Value = 0.364715860899067
$ws.Range("K19").Value = 0.08396946564885496
$ws.Range("M19").Value = 0.02120441051738762
$ws.Range("N19").Value = 0.0008481764206955047
$ws.Range("O19").Value = 0.07888040712468193
$ws.Range("S19").Value = 0.1068702290076336
